$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original plain-text storage.
# Excel auto-converts numeric-looking strings (e.g. "0.9996") to real
# numbers when assigned via .Value, which would corrupt values like
# "0.07100" (trailing zero) into 0.071. Pre-formatting the column as
# Text before writing keeps every price an exact string, matching the
# source data (which stores these as inline strings, not numbers).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.279.00"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3
$ws.Range("D3").Value = "1.929.74"
$ws.Range("E3").Value = "  -0.68%  "

# Row 4
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "0.7463"
$ws.Range("E5").Value = "  +3.06%  "

# Row 6
$ws.Range("D6").Value = "249.66"
$ws.Range("E6").Value = "  -0.98%  "

# Row 7
$ws.Range("D7").Value = "0.9989"
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").Value = "0.3232"
$ws.Range("E8").Value = "  -3.38%  "

# Row 9
$ws.Range("D9").Value = "27.91"
$ws.Range("E9").Value = "  -3.21%  "

# Row 10
$ws.Range("D10").Value = "0.07100"
$ws.Range("E10").Value = "  -4.75%  "

# Row 11
$ws.Range("D11").Value = "0.7890"
$ws.Range("E11").Value = "  -3.86%  "

# Row 12
$ws.Range("D12").Value = "0.08027"
$ws.Range("E12").Value = "  -1.35%  "

# Row 13
$ws.Range("D13").Value = "1.930.51"
$ws.Range("E13").Value = "  -0.53%  "

# Row 14
$ws.Range("D14").Value = "5.391"
$ws.Range("E14").Value = "  -1.99%  "

# Row 15
$ws.Range("D15").Value = "94.73"
$ws.Range("E15").Value = "  -0.75%  "

# Row 16
$ws.Range("E16").Value = "  -2.06%  "

# Row 17
$ws.Range("D17").Value = "30.276.78"
$ws.Range("E17").Value = "  -0.32%  "

# Row 18
$ws.Range("D18").Value = "254.66"
$ws.Range("E18").Value = "  +0.02%  "

# Row 19
$ws.Range("D19").Value = "0.000008054"
$ws.Range("E19").Value = "  -3.88%  "

# Row 20
$ws.Range("D20").Value = "5.750"
$ws.Range("E20").Value = "  -2.71%  "

# Row 21
$ws.Range("D21").Value = "2.184.08"
$ws.Range("E21").Value = "  -0.56%  "

# Row 22
$ws.Range("D22").Value = "0.9988"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("D23").Value = "0.9988"
$ws.Range("E23").Value = "  -0.06%  "

# Row 24
$ws.Range("D24").Value = "6.838"
$ws.Range("E24").Value = "  -2.25%  "

# Row 25
$ws.Range("D25").Value = "9.583"
$ws.Range("E25").Value = "  -3.35%  "

# Row 26
$ws.Range("D26").Value = "163.70"
$ws.Range("E26").Value = "  +0.96%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "19.11"
$ws.Range("E27").Value = "  -1.61%  "

# Row 28
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "0.1341"
$ws.Range("E28").Value = "  +1.26%  "

# Row 29
$ws.Range("D29").Value = "2.307"
$ws.Range("E29").Value = "  -4.88%  "

# Row 30
$ws.Range("D30").Value = "1.357"
$ws.Range("E30").Value = "  +0.49%  "

# Row 31
$ws.Range("D31").Value = "1.534"
$ws.Range("E31").Value = "  -2.50%  "

# Row 32
$ws.Range("D32").Value = "4.431"
$ws.Range("E32").Value = "  -0.95%  "

# Row 33
$ws.Range("D33").Value = "4.152"
$ws.Range("E33").Value = "  -2.76%  "

# Row 34
$ws.Range("D34").Value = "0.05114"
$ws.Range("E34").Value = "  -3.88%  "

# Row 35
$ws.Range("D35").Value = "1.293"
$ws.Range("E35").Value = "  -1.36%  "

# Row 36
$ws.Range("D36").Value = "0.7490"
$ws.Range("E36").Value = "  -1.76%  "

# Row 37
$ws.Range("D37").Value = "2.768"
$ws.Range("E37").Value = "  +0.65%  "

# Row 38
$ws.Range("D38").Value = "0.01978"
$ws.Range("E38").Value = "  -0.89%  "

# Row 39
$ws.Range("D39").Value = "2.800"
$ws.Range("E39").Value = "  -1.78%  "

# Row 40
$ws.Range("D40").Value = "78.24"
$ws.Range("E40").Value = "  -3.96%  "

# Row 41
$ws.Range("D41").Value = "6.407"
$ws.Range("E41").Value = "  -3.13%  "

# Row 42
$ws.Range("D42").Value = "0.4518"
$ws.Range("E42").Value = "  -1.45%  "

# Row 43
$ws.Range("D43").Value = "1.990"
$ws.Range("E43").Value = "  -2.94%  "

# Row 44
$ws.Range("D44").Value = "0.8426"
$ws.Range("E44").Value = "  -0.30%  "

# Row 45
$ws.Range("D45").Value = "0.9988"
$ws.Range("E45").Value = "  -0.12%  "

# Row 46
$ws.Range("D46").Value = "101.65"
$ws.Range("E46").Value = "  -1.53%  "

# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.579"
$ws.Range("E47").Value = "  +0.61%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.761"
$ws.Range("E48").Value = "  -0.77%  "

# Row 49
$ws.Range("D49").Value = "977.16"
$ws.Range("E49").Value = "  +11.02%  "

# Row 50
$ws.Range("D50").Value = "36.86"
$ws.Range("E50").Value = "  -0.68%  "

# Row 51
$ws.Range("D51").Value = "0.4199"
$ws.Range("E51").Value = "  -0.55%  "
